# "Generate Report for Handback"
#
# The handback report was regenerated for the first data row (row 2) on
# both locale sheets, producing fresh "Correspond Handoff Datetime" (col D)
# and "Correspond Handback DateTime" (col G) timestamps. Row 3 is untouched.
$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-22 16:40:45"
$wsZh.Range("G2").Value = "2016-02-22 16:41:33"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-22 16:40:56"
$wsDe.Range("G2").Value = "2016-02-22 16:41:56"
